$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values per row (columns D, L, M, N, O, P, R, S)
$data = @{
    2 = @{ D = 44544; L = "Segunda"; M = 250; N = 20000; O = 22000; P = 21000; R = "Provincia de San Felipe de Aconcagua"; S = 1167 }
    3 = @{ D = 44169; L = "Primera"; M = 250; N = 20000; O = 22000; P = 21000; R = "Provincia de San Felipe de Aconcagua"; S = 1167 }
    4 = @{ D = 44174; L = "Primera"; M = 300; N = 19000; O = 20000; P = 19500; R = "Región Metropolitana"; S = 1083 }
    5 = @{ D = 44524; L = "Segunda"; M = 200; N = 27000; O = 28000; P = 27500; R = "Provincia de San Felipe de Aconcagua"; S = 1528 }
    6 = @{ D = 44160; L = "Primera"; M = 250; N = 24000; O = 25000; P = 24500; R = "Provincia de San Felipe de Aconcagua"; S = 1361 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("L$row").Value = $vals.L
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("N$row").Value = $vals.N
    $ws.Range("O$row").Value = $vals.O
    $ws.Range("P$row").Value = $vals.P
    $ws.Range("R$row").Value = $vals.R
    $ws.Range("S$row").Value = $vals.S
}
